# Commit: "I sent the wrong saved draft"
# The prior save swapped the diff formulas (B-C instead of C-B, etc.) for the
# FY17/18/19 budget-vs-actual comparison columns. This restores the intended
# sign convention (Actual - Budget) for columns D, I and N on metro_budget,
# which ripples into the dependent % and RANK columns (E/F, J/K, O/P) and all
# of the downstream VLOOKUP/XLOOKUP/INDEX example tables.
# It also repositions the "Department Results" chart and updates the saved
# scroll position of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# --- Fix the sign of the budget/actual delta formulas -----------------
# Row 3 holds the shared-formula masters for D3:D52, I3:I52 and N3:N52 (and
# the dependent E/F, J/K, O/P columns recompute automatically), but row 2 is
# not part of those shared groups and needs its own formulas updated too.

$ws.Range("D2").Formula = "=(C2-B2)"
$ws.Range("I2").Formula = "=(H2-G2)"
$ws.Range("N2").Formula = "=M2-L2"

$ws.Range("D3:D52").Formula = "=(C3-B3)"
$ws.Range("I3:I52").Formula = "=(H3-G3)"
$ws.Range("N3:N52").Formula = "=M3-L3"

# --- Reposition / resize the "Department Results" chart ---------------
$co = $ws.ChartObjects(1)
$co.Left = 604.9100775098425
$co.Top = 1003.0
$co.Width = 544.337890625
$co.Height = 274.5
